# Staging.StrategicElement.xlsx was resaved with the BusinessKey / StrategicElement_ID
# header labels swapped: the "StrategicElement_ID" header moved from column A to column D,
# and "BusinessKey" moved from column D to column A.
#
# Apply that swap directly on the two header cells in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "BusinessKey"
$ws.Range("D2").Value = "StrategicElement_ID"
